# Apply updated experiment result values (no normalization) for rows 2-13
# on columns B, C, E, F, H, J, K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are given per-row as: B, C, E, F, H, J, K
$data = @{
    2  = @{ B=5; C=4; E=5; F=3; H=1; J=9;  K=374 }
    3  = @{ B=5; C=5; E=5; F=5; H=2; J=2;  K=12 }
    4  = @{ B=2; C=1; E=2; F=2; H=1; J=9;  K=374 }
    5  = @{ B=1; C=1; E=2; F=2; H=2; J=2;  K=12 }
    6  = @{ B=2; C=2; E=4; F=4; H=2; J=3;  K=22 }
    7  = @{ B=2; C=2; E=2; F=2; H=1; J=9;  K=374 }
    8  = @{ B=3; C=3; E=5; F=5; H=5; J=3;  K=22 }
    9  = @{ B=1; C=1; E=2; F=2; H=2; J=3;  K=22 }
    10 = @{ B=1; C=1; E=1; F=1; H=1; J=3;  K=22 }
    11 = @{ B=1; C=1; E=2; F=2; H=2; J=3;  K=22 }
    12 = @{ B=3; C=3; E=7; F=7; H=3; J=3;  K=22 }
    13 = @{ B=3; C=3; E=2; F=2; H=7; J=10; K=380 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
}
